# Apply weekly fruit/vegetable price data shuffle to the "Caqui" sheet.
# The D (Fecha), M (Volumen), Q (Unidad de comercialización), S (Precio $/Kg)
# and T (Kg / unidad) values for rows 2-10 are re-permuted across rows while
# all other columns remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current (pre-edit) values for the columns that get shuffled,
# keyed by row number, before any values are overwritten.
$cols = @("D", "M", "Q", "S", "T")
$original = @{}
foreach ($r in 2..10) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $original[$r] = $rowVals
}

# Destination row -> source row (the row whose original D/M/Q/S/T values
# should end up in the destination row).
$mapping = @{
    2  = 6
    3  = 7
    4  = 3
    5  = 10
    6  = 5
    7  = 9
    8  = 4
    9  = 8
    10 = 2
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $src = $original[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $src[$c]
    }
}
